# Horarios actualizados Linea 141 - 567
#
# A new scrape snapshot (taken at 07:49:32) is merged into the previous one
# (07:38:09) across the three worksheets: LP1912, LP1912-215, 6203-6173.
# Each new reading is inserted into its correct chronological
# (Hora_Llegada) slot, pushing the rows that were already below it one (or
# more) positions further down - exactly like the source diff shows.
#
# Insertion rows below are expressed in the ORIGINAL (pre-edit) row
# numbering and are applied from the bottom of the sheet upwards, so an
# earlier (smaller-numbered) insertion point is never disturbed by a later
# one that is still queued above it.

$wb = $excel.ActiveWorkbook

function Set-Row($ws, $r, $a, $b, $c, $d, $e) {
    $ws.Range("A$r").Value = $a
    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = $c
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
}

function Insert-Rows($ws, $atRow, $count) {
    $last = $atRow + $count - 1
    $ws.Range("A${atRow}:A${last}").EntireRow.Insert()
}

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 07:49:32"
$ws1.Range("A3").Value = "Total filas: 51"

# Work bottom-up using the ORIGINAL row numbers.
Insert-Rows $ws1 49 3
Set-Row $ws1 49 "07:49:32" "09:39" "15_ABASTO" 110 "LP1912"
Set-Row $ws1 50 "07:49:32" "09:41" "11_ETCHEVERRY" 112 "LP1912"
Set-Row $ws1 51 "07:49:32" "09:43" "16_P MOR-SANTA ANA" 114 "LP1912"

Insert-Rows $ws1 46 1
Set-Row $ws1 46 "07:49:32" "09:17" "14_ABASTO" 88 "LP1912"

Insert-Rows $ws1 45 1
Set-Row $ws1 45 "07:49:32" "09:12" "27_EL RETIRO" 83 "LP1912"

Insert-Rows $ws1 42 1
Set-Row $ws1 42 "07:49:32" "08:53" "215B_EL PATO" 64 "LP1912"

Insert-Rows $ws1 40 1
Set-Row $ws1 40 "07:49:32" "08:47" "23_HERNANDEZ" 58 "LP1912"

Insert-Rows $ws1 29 1
Set-Row $ws1 29 "07:49:32" "07:49" "215A_EL PATO" 0 "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 07:49:32"
$ws2.Range("A3").Value = "Total filas: 13"

Insert-Rows $ws2 15 1
Set-Row $ws2 15 "07:49:32" "08:53" "215B_EL PATO" 64 "LP1912"

Insert-Rows $ws2 13 1
Set-Row $ws2 13 "07:49:32" "07:49" "215A_EL PATO" 0 "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 07:49:32"
$ws3.Range("A3").Value = "Total filas: 7"

# New row appended at the end (row 12 did not exist before, so no shift is
# needed - it simply becomes a new last row).
Set-Row $ws3 12 "07:49:32" "09:21" "215A_LA PLATA" 92 "L6173"
